$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new prediction row (row 5) below the existing data.
$ws.Range("A5").Value = "22-03-2025"
$ws.Range("B5").Value = "Kolkata Knight Riders vs Royal Challengers Bengaluru"

# Column C (Toss) is unknown yet for this match, same as the other rows
# that have a blank-but-present Toss/Winner cell (e.g. E2, E3, C4).
# Copying one of those empty cells keeps C5 present-but-empty instead of
# leaving the cell missing entirely.
$ws.Range("E2").Copy($ws.Range("C5"))

$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Kolkata Knight Riders"
